$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 19 from
# serial date 45208 (2023-10-09) to 45212 (2023-10-13), preserving the
# existing date formatting on those cells.
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
